# Refresh weekly forecast data by one week:
#  - "Forecast Comparison" sheet: Week_Start_Date (col B) and forecast
#    numbers (cols D-H) for rows 2..17 are replaced with the newly
#    generated forecast (one week later than before).
#  - "Summary" sheet: derived statistics are updated to match.
#
# NOTE: Week_Start_Date and all Summary "Value" cells are stored as plain
# text in the workbook (not real dates/numbers). A leading apostrophe is
# used on assignment so Excel keeps them as text instead of auto-coercing
# date-looking / number-looking strings into date serials or numbers.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# New values for rows 2..17: Week_Start_Date (B), MyForecast (D),
# Amazon Mean (E), Amazon P70 (F), Amazon P80 (G), Amazon P90 (H)
$rows = @(
    @{ Row = 2;  Date = "2025-02-02"; D = 109; E = 122; F = 145; G = 168; H = 202 },
    @{ Row = 3;  Date = "2025-02-09"; D = 91;  E = 86;  F = 103; G = 120; H = 146 },
    @{ Row = 4;  Date = "2025-02-16"; D = 102; E = 85;  F = 102; G = 119; H = 145 },
    @{ Row = 5;  Date = "2025-02-23"; D = 109; E = 84;  F = 100; G = 117; H = 143 },
    @{ Row = 6;  Date = "2025-03-02"; D = 109; E = 84;  F = 101; G = 119; H = 149 },
    @{ Row = 7;  Date = "2025-03-09"; D = 108; E = 83;  F = 100; G = 118; H = 146 },
    @{ Row = 8;  Date = "2025-03-16"; D = 100; E = 77;  F = 93;  G = 111; H = 140 },
    @{ Row = 9;  Date = "2025-03-23"; D = 95;  E = 78;  F = 95;  G = 115; H = 146 },
    @{ Row = 10; Date = "2025-03-30"; D = 80;  E = 81;  F = 97;  G = 115; H = 143 },
    @{ Row = 11; Date = "2025-04-06"; D = 90;  E = 73;  F = 89;  G = 108; H = 137 },
    @{ Row = 12; Date = "2025-04-13"; D = 86;  E = 66;  F = 80;  G = 98;  H = 127 },
    @{ Row = 13; Date = "2025-04-20"; D = 86;  E = 66;  F = 81;  G = 99;  H = 128 },
    @{ Row = 14; Date = "2025-04-27"; D = 84;  E = 65;  F = 79;  G = 96;  H = 122 },
    @{ Row = 15; Date = "2025-05-04"; D = 79;  E = 61;  F = 74;  G = 92;  H = 121 },
    @{ Row = 16; Date = "2025-05-11"; D = 79;  E = 61;  F = 75;  G = 92;  H = 120 },
    @{ Row = 17; Date = "2025-05-18"; D = 75;  E = 58;  F = 71;  G = 88;  H = 117 }
)

foreach ($r in $rows) {
    $wsForecast.Cells.Item($r.Row, 2).Value = "'" + $r.Date
    $wsForecast.Cells.Item($r.Row, 4).Value = $r.D
    $wsForecast.Cells.Item($r.Row, 5).Value = $r.E
    $wsForecast.Cells.Item($r.Row, 6).Value = $r.F
    $wsForecast.Cells.Item($r.Row, 7).Value = $r.G
    $wsForecast.Cells.Item($r.Row, 8).Value = $r.H
}

# Update Summary sheet values (all stored as text)
$wsSummary.Cells.Item(2, 2).Value  = "'2022-12-25 to 2025-01-26"
$wsSummary.Cells.Item(4, 2).Value  = "'212"
$wsSummary.Cells.Item(5, 2).Value  = "'81"
$wsSummary.Cells.Item(6, 2).Value  = "'66"
$wsSummary.Cells.Item(7, 2).Value  = "'62"
$wsSummary.Cells.Item(8, 2).Value  = "'8787 units"
$wsSummary.Cells.Item(9, 2).Value  = "'1484"
$wsSummary.Cells.Item(10, 2).Value = "'823"
$wsSummary.Cells.Item(11, 2).Value = "'411"
$wsSummary.Cells.Item(12, 2).Value = "'109"
$wsSummary.Cells.Item(13, 2).Value = "'2025-02-23"
$wsSummary.Cells.Item(14, 2).Value = "'75"
$wsSummary.Cells.Item(15, 2).Value = "'2025-05-18"
